# Append run: 2025-11-07 12:35 JST
# - Two brand-new listings are spliced into the existing ranking (at row 2
#   and at row 6), pushing the previously-scraped rows down.
# - Every row's "取得日時" (fetched-at) timestamp is refreshed to the new run.
# - The two rows that land at the bottom of the sheet (old rows 9 and 10)
#   need fresh Hyperlink objects on their URL cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Make room for the two new listings ------------------------------
# Inserting shifts existing rows/styles down (Excel-style), but leaves the
# worksheet's existing Hyperlinks collection entries pointing at their old
# `ref`s - matching upstream's own scraper-export behaviour.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(6).Insert()

# --- Row 2 (new): AI-automated secondhand-brand resale tooling --------
$ws.Range("A2").Value = "2025-11-07 12:35:55"
$ws.Range("B2").Value = "中古ブランド品リサーチとEC出品作業をAIで自動化するツール開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5429252"
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("G2").Value = 480
$ws.Range("H2").Value = "🔥AI,Ai ◆ツール,開発"

# --- Row 6 (new): mobile-app test-ops contract -------------------------
$ws.Range("A6").Value = "2025-11-07 12:35:55"
$ws.Range("B6").Value = "【急募】モバイルアプリ テスト業務 委託募集(3 - 4週間)"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5429220"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("G6").Value = 38
$ws.Range("H6").Value = "◇アプリ"

# --- Refresh the scrape timestamp on every data row ---------------------
$ws.Range("A2:A12").Value = "2025-11-07 12:35:55"

# --- New hyperlink objects for the rows that shifted to the tail --------
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5428970")
$ws.Range("F11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5428509")
$ws.Range("F12").Style = "Hyperlink"
